# Apply the changes described by the diff:
#  1. Split the run "u'il " (curly apostrophe) into "u'il" followed by a new
#     inline <comment>c_105v_03</comment> marker run-set and a trailing space run.
#  2. Add a footer distance (w:footer="720") to the section's page margins.

$d = $word.ActiveDocument

# --- Part 1: insert the <comment>c_105v_03</comment> marker -----------------

$apos = [char]0x2019

# "u'il " (with the curly apostrophe) occurs twice in the document ("qu'il se
# ferme" and "qu'il n'y demeure"); disambiguate by searching for the longer,
# unique phrase "u'il n'y" and then narrow down to just the first five
# characters ("u'il ") of that match.
$probe = $d.Content
$found = $probe.Find.Execute("u" + $apos + "il n" + $apos + "y", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate target text 'u$($apos)il n$($apos)y'"
}

$start = $probe.Start
$target = $d.Range($start, $start + 5)
if ($target.Text -ne ("u" + $apos + "il ")) {
    throw "Unexpected target text: [$($target.Text)]"
}

# Shrink the run's text so it no longer includes the trailing space.
$target.Text = "u" + $apos + "il"

# Insert the new runs right after "u'il", each with its own formatting,
# mirroring the existing "<ms><bp>...</bp></ms>" markup runs already present
# in the document.
$insPos = $start + 4
$r1 = $d.Range($insPos, $insPos)
$r1.InsertAfter("<comment>")
$r1.Font.Name = "Courier New"
$r1.Font.Color = 16711680
$r1.Font.Size = 9

$insPos = $insPos + 9
$r2 = $d.Range($insPos, $insPos)
$r2.InsertAfter("c_105v_03")
$r2.Font.Size = 8

$insPos = $insPos + 9
$r3 = $d.Range($insPos, $insPos)
$r3.InsertAfter("</comment>")
$r3.Font.Name = "Courier New"
$r3.Font.Color = 16711680
$r3.Font.Size = 9

$insPos = $insPos + 10
$r4 = $d.Range($insPos, $insPos)
$r4.InsertAfter(" ")

# --- Part 2: set the footer distance on the section page margins -----------

foreach ($sec in $d.Sections) {
    $sec.PageSetup.FooterDistance = 36
}
